# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (and so to every
#                             slide) -- currently the "Integral" palette.
#   ppt/theme/theme2.xml  -> bound to the notes master -- currently the
#                             stock "Office Theme" palette.
# The authored commit swaps the two palettes: the slide master (and thus
# the whole deck) switches to the default "Office Theme" 12-colour scheme,
# while the notes master keeps/receives the old "Integral" colours.
#
# PowerPoint's theme-colour object model is reached from the slide master
# via Theme.ThemeColorScheme (12 entries, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Each entry's .RGB is a
# standard Windows COLORREF (0xBBGGRR, i.e. R + G*256 + B*65536) -- not a
# plain RRGGBB hex read -- so every literal below is that conversion of
# the target srgbClr hex value.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
